# Multiply every "Total" (column AB) value in the data rows (2-524) by 1.5
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 524
$col = 28  # Column AB

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = $val * 1.5
    }
}
